$d = $word.ActiveDocument
$t = $d.Tables.Item(1)

function Set-CellText($rowIndex, $text) {
    $cell = $t.Rows.Item($rowIndex).Cells.Item(1)
    $cell.Range.Text = $text
}

# Simple single-value cell replacements (1-based row indices)
Set-CellText 1 "0M"
Set-CellText 2 "0M"
Set-CellText 3 "0M"
Set-CellText 4 "400"
Set-CellText 6 "0.00052"
Set-CellText 7 "0.00021"
Set-CellText 8 "0.00006"
Set-CellText 9 "0.00034"
Set-CellText 10 "0.00040"
Set-CellText 11 "0.00045"
Set-CellText 12 "0.08421"

# Rows 44-46 (1-based) had multiple runs separated by tabs; collapse to a single value
Set-CellText 44 "99.96"
Set-CellText 45 "0.08"
Set-CellText 46 "215"
